$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 182.4
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 185.09091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2357.6
$ws.Range("J43").Value = 2897
$ws.Range("L43").Value = 2897
$ws.Range("N43").Value = -3035

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2718.625
$ws.Range("I51").Value = 1877
$ws.Range("J51").Value = 2999.1667
$ws.Range("K51").Value = 1877
$ws.Range("L51").Value = 2999.1667
$ws.Range("M51").Value = -1393
$ws.Range("N51").Value = -3967.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 781.8333
$ws.Range("I104").Value = 781.8333
$ws.Range("K104").Value = 2345.4999
$ws.Range("M104").Value = -598.4998999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4747
$ws.Range("I116").Value = 4766
$ws.Range("J116").Value = 4690
$ws.Range("K116").Value = 4766
$ws.Range("L116").Value = 4690
$ws.Range("M116").Value = -1324
$ws.Range("N116").Value = -11574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1190.3334
$ws.Range("I32").Value = 1190.3334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1190.3334
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -903.3334
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1954.1538
$ws.Range("I61").Value = 1440
$ws.Range("J61").Value = 3668
$ws.Range("K61").Value = 1440
$ws.Range("L61").Value = 3668
$ws.Range("M61").Value = -1228
$ws.Range("N61").Value = -4092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2826.3
$ws.Range("I132").Value = 2826.3
$ws.Range("K132").Value = 8478.900000000001
$ws.Range("M132").Value = -5948.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1954.1538
$ws.Range("I136").Value = 1440
$ws.Range("J136").Value = 3668
$ws.Range("K136").Value = 4320
$ws.Range("L136").Value = 11004
$ws.Range("M136").Value = -1770
$ws.Range("N136").Value = -16104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 915
$ws.Range("I20").Value = 789.5
$ws.Range("J20").Value = 998.6667
$ws.Range("K20").Value = 789.5
$ws.Range("L20").Value = 998.6667
$ws.Range("M20").Value = -542.5
$ws.Range("N20").Value = -1492.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4517.1816
$ws.Range("I86").Value = 1157.8
$ws.Range("J86").Value = 7316.6665
$ws.Range("K86").Value = 1157.8
$ws.Range("L86").Value = 7316.6665
$ws.Range("M86").Value = -34.79999999999995
$ws.Range("N86").Value = -9562.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4517.1816
$ws.Range("I89").Value = 1157.8
$ws.Range("J89").Value = 7316.6665
$ws.Range("K89").Value = 5789
$ws.Range("L89").Value = 36583.3325
$ws.Range("M89").Value = -173
$ws.Range("N89").Value = -47815.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2434.4546
$ws.Range("I99").Value = 2434.4546
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2434.4546
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -936.4546
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2434.4546
$ws.Range("I126").Value = 2434.4546
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7303.3638
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4833.3638
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4435.7617
$ws.Range("J80").Value = 5787.75
$ws.Range("L80").Value = 17363.25
$ws.Range("N80").Value = -19235.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4435.7617
$ws.Range("J83").Value = 5787.75
$ws.Range("L83").Value = 52089.75
$ws.Range("N83").Value = -61449.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1000
$ws.Range("I109").Value = 1000
$ws.Range("K109").Value = 3000
$ws.Range("M109").Value = -1960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 806.875
$ws.Range("J97").Value = 953
$ws.Range("L97").Value = 953
$ws.Range("N97").Value = -1945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3096.8333
$ws.Range("I132").Value = 2816.4
$ws.Range("K132").Value = 8449.200000000001
$ws.Range("M132").Value = -5919.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 95331.664
$ws.Range("J134").Value = 95331.664
$ws.Range("L134").Value = 285994.992
$ws.Range("N134").Value = -291064.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1633.3334
$ws.Range("I16").Value = 1800
$ws.Range("J16").Value = 1550
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 1550
$ws.Range("M16").Value = -1630
$ws.Range("N16").Value = -1890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 731.0417
$ws.Range("I55").Value = 789.3333
$ws.Range("J55").Value = 672.75
$ws.Range("K55").Value = 789.3333
$ws.Range("L55").Value = 672.75
$ws.Range("M55").Value = -616.3333
$ws.Range("N55").Value = -1018.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 37087
$ws.Range("J41").Value = 37823
$ws.Range("L41").Value = 37823
$ws.Range("N41").Value = -38603

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 31838.666
$ws.Range("J69").Value = 33135
$ws.Range("L69").Value = 33135
$ws.Range("N69").Value = -34633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 31838.666
$ws.Range("J72").Value = 33135
$ws.Range("L72").Value = 99405
$ws.Range("N72").Value = -106893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6567.6772
$ws.Range("I126").Value = 5843.9375
$ws.Range("J126").Value = 7339.6665
$ws.Range("K126").Value = 17531.8125
$ws.Range("L126").Value = 22018.9995
$ws.Range("M126").Value = -15061.8125
$ws.Range("N126").Value = -26958.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3585.4375
$ws.Range("I136").Value = 2901.818
$ws.Range("K136").Value = 8705.454000000002
$ws.Range("M136").Value = -6155.454000000002
